$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 (shifts old rows 48-62 down to 49-63)
$ws.Rows.Item(48).Insert()

$ws.Cells.Item(48, 1).Value = 11
$ws.Cells.Item(48, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value = "Bíobío"
$ws.Cells.Item(48, 4).Value = 45007
$ws.Cells.Item(48, 5).Value = 8
$ws.Cells.Item(48, 6).Value = 100112030
$ws.Cells.Item(48, 7).Value = "Poroto granado"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 100
$ws.Cells.Item(48, 11).Value = 30000
$ws.Cells.Item(48, 12).Value = 32000
$ws.Cells.Item(48, 13).Value = 31000
$ws.Cells.Item(48, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value = 1240
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
